$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header row: USN, Name, Age, Phone
$ws.Range("A1").Value = "USN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"

# Leave the active selection on D1, matching the authored workbook state
$ws.Range("D1").Select() | Out-Null
